# Automatic update of files.
# Column C ("Förändrad") holds a date serial value that is bumped by one day
# (45779 -> 45780) for every data row (rows 2 through 43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 43 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45779) {
        $cell.Value2 = 45780
    }
}
